# Generate Report for Handback
# Update the timestamp strings recorded in the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 03:04:46"

# "zh-cn" sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 03:04:41"
$wsZhCn.Range("K2").Value = "2016-08-25 03:04:59"

# "de-de" sheet: Latest HO Xliff Generate Date is shared with Overview (H2),
# Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 03:04:46"
$wsDeDe.Range("K2").Value = "2016-08-25 03:05:14"
